$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each entry: row, column, new text
$updates = @(
    @(1, 1, "44÷3=14, 2"),
    @(1, 2, "10÷3=3, 1"),
    @(1, 3, "22÷3=7, 1"),
    @(1, 4, "28÷5=5, 3"),
    @(1, 5, "55÷5=11, 0"),

    @(5, 1, "45÷5=9, 0"),
    @(5, 2, "66÷6=11, 0"),
    @(5, 3, "72÷7=10, 2"),
    @(5, 4, "95÷7=13, 4"),
    @(5, 5, "51÷4=12, 3"),

    @(9, 1, "29÷4=7, 1"),
    @(9, 2, "83÷6=13, 5"),
    @(9, 3, "70÷7=10, 0"),
    @(9, 4, "81÷5=16, 1"),
    @(9, 5, "11÷2=5, 1"),

    @(13, 1, "90÷4=22, 2"),
    @(13, 2, "29÷6=4, 5"),
    @(13, 3, "16÷6=2, 4"),
    @(13, 4, "46÷8=5, 6"),
    @(13, 5, "20÷2=10, 0"),

    @(17, 1, "90÷6=15, 0"),
    @(17, 2, "28÷8=3, 4"),
    @(17, 3, "43÷4=10, 3"),
    @(17, 4, "94÷2=47, 0"),
    @(17, 5, "89÷8=11, 1")
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $newText = $u[2]
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $newText
}
